# Silverlake.py results — widen the existing index/value table with two
# extra "Unnamed: 0.x" index columns and append ten more rows of data.
#
# Original layout (rows 1-11):
#   A = index (styled), B = "Unnamed: 0", C = "   SILVER_FOR", D = "SILVER_FOR"
# New layout (rows 1-21):
#   A = index (styled), B = "Unnamed: 0.2", C = "Unnamed: 0.1", D = "Unnamed: 0",
#   E = "   SILVER_FOR", F = "SILVER_FOR"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: shift old columns B:D (11 rows incl. header) two columns right
# to D:F. A single block copy carries over values, header text, number
# formatting/styles and the existing blank placeholder cells in one shot.
$ws.Range("B1:D11").Copy($ws.Range("D1:F11"))

# --- Phase 2: new header cells B1/C1 — reuse the header style that now
# lives on D1 (bold/bordered/centered), then set their text.
$ws.Range("D1").Copy($ws.Range("B1:C1"))
$ws.Range("B1").Value = "Unnamed: 0.2"
$ws.Range("C1").Value = "Unnamed: 0.1"

# --- Phase 3: extend column A (styled index) down through row 21.
$ws.Range("A2:A11").Copy($ws.Range("A12:A21"))
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $i
}

# --- Phase 4: new column B ("Unnamed: 0.2") data, rows 2-16 -> 0..14
for ($i = 0; $i -lt 15; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $i
}

# --- Phase 5: new column C ("Unnamed: 0.1") data, rows 2-11 -> 0..9
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $i
}

# --- Phase 6: blank staircase placeholders for the newly-added rows,
# reusing an existing blank cell (D2, an empty "SILVER_FOR" slot) as the
# copy source so the blank cells come through the same way the sheet's
# original blanks do. ClearContents first works around the engine
# quirk where Copy-ing onto the top-left cell of a destination block
# silently no-ops if that cell already holds a value.
$ws.Range("C12:E16").ClearContents()
$ws.Range("D2").Copy($ws.Range("C12:E16"))
$ws.Range("B17:E21").ClearContents()
$ws.Range("D2").Copy($ws.Range("B17:E21"))

# --- Phase 7: new rows 12-21, column F ("SILVER_FOR") values.
$ws.Cells.Item(12, 6).Value = 31.76244298992617
$ws.Cells.Item(13, 6).Value = 32.23934168175691
$ws.Cells.Item(14, 6).Value = 32.5524573182023
$ws.Cells.Item(15, 6).Value = 32.38868114596181
$ws.Cells.Item(16, 6).Value = 32.10415671664822
$ws.Cells.Item(17, 6).Value = 28.81685345771996
$ws.Cells.Item(18, 6).Value = 28.51061683625062
$ws.Cells.Item(19, 6).Value = 27.67123246313037
$ws.Cells.Item(20, 6).Value = 27.34084839285276
$ws.Cells.Item(21, 6).Value = 26.94339908891732

Write-Host "Applied Silverlake.py results update"
